$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 145.84616
$ws.Range("I2").Value = 145.84616
$ws.Range("K2").Value = 145.84616
$ws.Range("M2").Value = -32.84616
$ws.Range("H33").Value = 370.91177
$ws.Range("I33").Value = 214.23077
$ws.Range("J33").Value = 880.125
$ws.Range("K33").Value = 214.23077
$ws.Range("L33").Value = 880.125
$ws.Range("M33").Value = 14.76922999999999
$ws.Range("N33").Value = -1338.125
$ws.Range("H55").Value = 97.23077000000001
$ws.Range("I55").Value = 97.14286
$ws.Range("J55").Value = 97.333336
$ws.Range("K55").Value = 97.14286
$ws.Range("L55").Value = 97.333336
$ws.Range("M55").Value = 116.85714
$ws.Range("N55").Value = -525.333336
$ws.Range("H62").Value = 8891.267
$ws.Range("I62").Value = 7521.1113
$ws.Range("K62").Value = 7521.1113
$ws.Range("M62").Value = -6897.1113
$ws.Range("H65").Value = 8891.267
$ws.Range("I65").Value = 7521.1113
$ws.Range("K65").Value = 37605.5565
$ws.Range("M65").Value = -34485.5565
$ws.Range("H74").Value = 5056.3
$ws.Range("I74").Value = 4030.5
$ws.Range("J74").Value = 7449.8335
$ws.Range("K74").Value = 4030.5
$ws.Range("L74").Value = 7449.8335
$ws.Range("M74").Value = -3094.5
$ws.Range("N74").Value = -9321.833500000001
$ws.Range("H77").Value = 5056.3
$ws.Range("I77").Value = 4030.5
$ws.Range("J77").Value = 7449.8335
$ws.Range("K77").Value = 20152.5
$ws.Range("L77").Value = 37249.1675
$ws.Range("M77").Value = -15472.5
$ws.Range("N77").Value = -46609.1675
$ws.Range("H100").Value = 3127.0527
$ws.Range("J100").Value = 4787.375
$ws.Range("L100").Value = 4787.375
$ws.Range("N100").Value = -5869.375
$ws.Range("H135").Value = 1667
$ws.Range("I135").Value = 1242.3125
$ws.Range("K135").Value = 11180.8125
$ws.Range("M135").Value = -8645.8125
$ws.Range("H138").Value = 2908.2322
$ws.Range("I138").Value = 2557.4119
$ws.Range("J138").Value = 3061.1538
$ws.Range("K138").Value = 7672.2357
$ws.Range("L138").Value = 9183.4614
$ws.Range("M138").Value = -2532.2357
$ws.Range("N138").Value = -19463.4614

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1588.25
$ws.Range("I2").Value = 1031.6923
$ws.Range("J2").Value = 4000
$ws.Range("K2").Value = 1031.6923
$ws.Range("L2").Value = 4000
$ws.Range("M2").Value = -918.6922999999999
$ws.Range("N2").Value = -4226
$ws.Range("H116").Value = 1588.25
$ws.Range("I116").Value = 1031.6923
$ws.Range("J116").Value = 4000
$ws.Range("K116").Value = 1031.6923
$ws.Range("L116").Value = 4000
$ws.Range("M116").Value = 1262.3077
$ws.Range("N116").Value = -8588
$ws.Range("H122").Value = 4007.5847
$ws.Range("I122").Value = 2969.2654
$ws.Range("K122").Value = 8907.796200000001
$ws.Range("M122").Value = -6457.796200000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1588.25
$ws.Range("I3").Value = 1031.6923
$ws.Range("J3").Value = 4000
$ws.Range("K3").Value = 1031.6923
$ws.Range("L3").Value = 4000
$ws.Range("M3").Value = -917.6922999999999
$ws.Range("N3").Value = -4228
$ws.Range("H105").Value = 3232.875
$ws.Range("I105").Value = 1662.8334
$ws.Range("K105").Value = 1662.8334
$ws.Range("M105").Value = 84.16660000000002
$ws.Range("H126").Value = 39339.5
$ws.Range("J126").Value = 39339.5
$ws.Range("L126").Value = 39339.5
$ws.Range("N126").Value = -49219.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2047.7826
$ws.Range("I132").Value = 2057.8948
$ws.Range("J132").Value = 1999.75
$ws.Range("K132").Value = 6173.6844
$ws.Range("L132").Value = 5999.25
$ws.Range("M132").Value = -3643.6844
$ws.Range("N132").Value = -11059.25
$ws.Range("H134").Value = 2278.5483
$ws.Range("I134").Value = 1676.04
$ws.Range("J134").Value = 4789
$ws.Range("K134").Value = 5028.12
$ws.Range("L134").Value = 14367
$ws.Range("M134").Value = -2493.12
$ws.Range("N134").Value = -19437

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 584.1177
$ws.Range("J38").Value = 769.6
$ws.Range("L38").Value = 2308.8
$ws.Range("N38").Value = -3002.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4514.9414
$ws.Range("I70").Value = 4627
$ws.Range("K70").Value = 4627
$ws.Range("M70").Value = -4357
$ws.Range("H73").Value = 4514.9414
$ws.Range("I73").Value = 4627
$ws.Range("K73").Value = 4627
$ws.Range("M73").Value = -3691
$ws.Range("H80").Value = 2226.5557
$ws.Range("J80").Value = 2999.8333
$ws.Range("L80").Value = 2999.8333
$ws.Range("N80").Value = -4995.8333
$ws.Range("H83").Value = 2226.5557
$ws.Range("J83").Value = 2999.8333
$ws.Range("L83").Value = 14999.1665
$ws.Range("N83").Value = -24983.1665
$ws.Range("H132").Value = 2722.389
$ws.Range("I132").Value = 2550.5
$ws.Range("K132").Value = 7651.5
$ws.Range("M132").Value = -5121.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6453.4585
$ws.Range("I40").Value = 6544.3
$ws.Range("J40").Value = 5999.25
$ws.Range("K40").Value = 6544.3
$ws.Range("L40").Value = 5999.25
$ws.Range("M40").Value = -6408.3
$ws.Range("N40").Value = -6271.25
$ws.Range("H61").Value = 1686.25
$ws.Range("I61").Value = 1641.4286
$ws.Range("K61").Value = 1641.4286
$ws.Range("M61").Value = -1439.4286
$ws.Range("H113").Value = 1686.25
$ws.Range("I113").Value = 1641.4286
$ws.Range("K113").Value = 1641.4286
$ws.Range("M113").Value = 528.5714
$ws.Range("H122").Value = 3626.4644
$ws.Range("I122").Value = 3733.9285
$ws.Range("K122").Value = 11201.7855
$ws.Range("M122").Value = -8751.7855
$ws.Range("H132").Value = 2265.838
$ws.Range("I132").Value = 1963.6562
$ws.Range("K132").Value = 5890.9686
$ws.Range("M132").Value = -3360.9686

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 30025
$ws.Range("I40").Value = 30025
$ws.Range("K40").Value = 30025
$ws.Range("M40").Value = -29876
$ws.Range("H52").Value = 44007.6
$ws.Range("J52").Value = 47499
$ws.Range("L52").Value = 47499
$ws.Range("N52").Value = -47951
$ws.Range("H69").Value = 15000
$ws.Range("J69").Value = 15000
$ws.Range("L69").Value = 15000
$ws.Range("N69").Value = -16498
$ws.Range("H72").Value = 15000
$ws.Range("J72").Value = 15000
$ws.Range("L72").Value = 45000
$ws.Range("N72").Value = -52488
$ws.Range("H74").Value = 38687
$ws.Range("J74").Value = 64124.5
$ws.Range("L74").Value = 64124.5
$ws.Range("N74").Value = -65996.5
$ws.Range("H77").Value = 38687
$ws.Range("J77").Value = 64124.5
$ws.Range("L77").Value = 192373.5
$ws.Range("N77").Value = -201733.5
$ws.Range("H126").Value = 6647.8887
$ws.Range("I126").Value = 6482.385
$ws.Range("K126").Value = 19447.155
$ws.Range("M126").Value = -16977.155
